$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (style s="1") from an existing header cell (AC1) onto
# the three new header cells, then set their text values.
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 62
    $ws.Cells.Item($row, 31).Value = 99
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Host "done"
